$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append the new row of data (row 5) that arrived at 2025-05-01T11:00:11.932Z
$ws.Range("A5").Value = "2025-05-01T11:00:11.932Z"
$ws.Range("B5").Value = "NRC"
$ws.Range("C5").Value = "C3"
$ws.Range("D5").Value = "أخرى"
$ws.Range("E5").Value = "الصمود"
$ws.Range("F5").Value = "احمد"

# G5 ("32") must be stored as text (matching the rest of the sheet's
# "numbers stored as text" convention), not auto-converted to a number.
$ws.Range("G5").NumberFormat = "@"
$ws.Range("G5").Value = "32"
$ws.Range("G5").Style = "Normal"

# H5 is an (empty) text cell in the source row.
$ws.Range("H5").NumberFormat = "@"
$ws.Range("H5").Value = ""
$ws.Range("H5").Style = "Normal"
